$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4026.447856934395,
    4003.936478312752,
    4002.722642646088,
    4002.722642646088,
    4002.722642646088,
    3996.998694298975,
    3951.525545262272,
    3857.687076954601,
    3857.687076954601,
    3857.687076954601,
    3857.687076954601
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
